# "created a start menu"
# Add a new "Status" column (G) to the jobs table and mark every job as "Active".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header for column G
$ws.Range("G1").Value = "Status"

# Every data row (2-13) now has a Status of "Active"
$ws.Range("G2:G13").Value = "Active"

# Leave the selection on the last edited cell (matches the saved cursor position)
$ws.Range("H13").Select() | Out-Null
